# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" between "2021-Q2" and "总计" with the
#   fund-holding detail for the new quarter.
# - Add a new summary row for "2022-Q1" at the top of the "总计" sheet's
#   data (pushing the existing "2021-Q2" row down).

$wb = $excel.ActiveWorkbook

$wsQ2       = $wb.Worksheets.Item(1)   # "2021-Q2"
$wsTotalOld = $wb.Worksheets.Item(2)   # "总计" (before the new sheet is inserted)

# Helper: write a numeric-looking value into a cell as LITERAL TEXT
# (leading/trailing zeros matter, e.g. "005075", "0.0000"), then strip the
# explicit style that setting NumberFormat="@" stamps on the cell by
# pasting the format from a never-touched ("no style") cell on top of it,
# so the cell ends up with no explicit style attribute at all.
function Set-TextValue($ws, $addr, $value, $blankAddr) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $ws.Range($blankAddr).Copy()
    $r.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating "总计" so that it
#    starts out with the same header/row styling ("总计" uses the same
#    style index for its header + first column as the new sheet needs),
#    positioned right after "2021-Q2".
# ---------------------------------------------------------------------
$wsTotalOld.Copy($null, $wsQ2)
$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "2022-Q1"

# Extend the header styling from D1 across E1:H1 (copy format only).
$wsNew.Range("D1").Copy()
$wsNew.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row text.
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Give row 3's "A" cell the same styling as row 2's (carried over from the
# copied "总计" sheet) before filling in values.
$wsNew.Range("A2").Copy()
$wsNew.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsNew.Range("A2").Value = 0
$wsNew.Range("A3").Value = 1

# Plain-text (non numeric-looking) cells -- safe to assign directly.
$wsNew.Range("C2").Value = "富国研究量化精选混合"
$wsNew.Range("C3").Value = "新华鑫弘灵活配置混合"

# Numeric-looking values that must be preserved verbatim as text.
Set-TextValue $wsNew "B2" "005075" "C2"
Set-TextValue $wsNew "D2" "3.03"   "C2"
Set-TextValue $wsNew "E2" "94.40"  "C2"
Set-TextValue $wsNew "F2" "1.74"   "C2"
Set-TextValue $wsNew "G2" "0.0527" "C2"

Set-TextValue $wsNew "B3" "003739" "C3"
Set-TextValue $wsNew "D3" "0.01"   "C3"
Set-TextValue $wsNew "E3" "42.81"  "C3"
Set-TextValue $wsNew "F3" "0.47"   "C3"
Set-TextValue $wsNew "G3" "0.0000" "C3"

# Genuinely-numeric cells.
$wsNew.Range("H2").Value = 10
$wsNew.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new summary row on top (row 2),
#    pushing "2021-Q2" down to row 3.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(3)   # "总计" after the insert above

$wsTotal.Rows.Item(2).Insert()

# New row 2 inherits row-1's styling by default; give column A the same
# style as the (now) row-3 "A" cell, and strip the accidental styling
# that Insert() propagates onto B2:D2.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTotal.Range("C3").Copy()
$wsTotal.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.05
$wsTotal.Range("A3").Value = 1
